# Fill in the "digikey" order quantities (column U, rows 7-32) of the BOM
# worksheet, as described in the commit:
#   "Added README, filled in BOM to provide order list for DIGIKEY"
#
# This drives the existing LOOKUP()/SUMIF() formulas already present in the
# sheet (columns V/W and the roll-up cells in row 34 and the summary rows
# 35-60), which Excel will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$quantities = @{
    7  = 2
    8  = 6
    9  = 2
    10 = 2
    11 = 2
    12 = 9
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 2
    27 = 1
    28 = 2
    29 = 1
    30 = 1
    31 = 1
    32 = 3
}

foreach ($row in $quantities.Keys) {
    $ws.Range("U$row").Value = $quantities[$row]
}

# Recalculate the workbook so the dependent formulas (V/W columns, the
# purchased-parts summaries in row 34, and the per-vendor pick lists in
# rows 35-60) refresh their cached values.
$excel.CalculateFullRebuild()

# Reflect where the user ended up working in the sheet (last scroll
# position / selection) when the file was saved.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 19
$ws.Range("T37").Select()
